# Added the characterization factors of plastic resin from the MariLCA group
# + correct a mistake that was dividing the S matrix for non-plastic flows by 1,000,000.
#
# All edits target the "Canada" sheet (the only sheet whose tabSelected view
# was active / whose data changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Canada")

# Row 2 - Polyethylene terephthalate (PET) resins
$ws.Range("F2").Value = 0.23
$ws.Range("J2").Value = 0.01

# Row 3 - Other thermoplastic polyester resins
$ws.Range("J3").Value = 0.01

# Row 4 - Low-density polyethylene (LDPE) resins
$ws.Range("E4").Value = 0.3

# Row 5 - Linear low-density polyethylene (LLDPE) resins
$ws.Range("E5").Value = 0.3

# Row 6 - High-density polyethylene (HDPE) resins
$ws.Range("E6").Value = 0.17

# Row 10 - Polyvinyl chloride (PVC) resins
$ws.Range("E10").Value = 0.03
$ws.Range("F10").Value = 0.02

# Row 11 - Polypropylene (PP) resins
$ws.Range("E11").Value = 0.2

# Row 19 - Vehicles
$ws.Range("J19").Value = 0.08

# Match the saved cursor/selection position on the Canada sheet
$ws.Range("A24").Select()
